$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.185.99'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.833.16'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('D4').Value = "'0.9995"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'241.97"
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').Value = "'0.6577"
$ws.Range('E6').Value = '  -2.77%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = "'0.07402"
$ws.Range('E8').Value = '  -0.74%  '
$ws.Range('D9').Value = "'0.2933"
$ws.Range('E9').Value = '  -1.59%  '
$ws.Range('D10').Value = "'22.85"
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('D11').Value = "'0.07768"
$ws.Range('E11').Value = '  +1.30%  '
$ws.Range('D12').Value = '1.843.31'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').Value = "'4.996"
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('D14').Value = "'0.6661"
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('D15').Value = "'82.81"
$ws.Range('E15').Value = '  -4.20%  '
$ws.Range('D16').Value = "'6.103"
$ws.Range('E16').Value = '  -0.96%  '
$ws.Range('D17').Value = "'0.000008404"
$ws.Range('E17').Value = '  +1.74%  '
$ws.Range('D18').Value = '29.167.41'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').Value = '2.061.33'
$ws.Range('E19').Value = '  -0.86%  '
$ws.Range('D20').Value = "'227.12"
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D21').Value = "'12.45"
$ws.Range('E21').Value = '  -0.32%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').Value = "'7.127"
$ws.Range('E23').Value = '  -2.94%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').Value = "'159.09"
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('D26').Value = "'8.606"
$ws.Range('E26').Value = '  -1.04%  '
$ws.Range('D27').Value = "'0.1388"
$ws.Range('E27').Value = '  -3.11%  '
$ws.Range('D28').Value = "'17.93"
$ws.Range('E28').Value = '  -0.47%  '
$ws.Range('D29').Value = "'1.518"
$ws.Range('E29').Value = '  +1.02%  '
$ws.Range('D30').Value = "'4.112"
$ws.Range('E30').Value = '  -3.19%  '
$ws.Range('D31').Value = "'4.044"
$ws.Range('E31').Value = '  -2.08%  '
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('D33').Value = "'0.05272"
$ws.Range('E33').Value = '  -2.84%  '
$ws.Range('D34').Value = "'1.860"
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('D35').Value = "'0.7404"
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('D36').Value = "'1.140"
$ws.Range('E36').Value = '  +1.06%  '
$ws.Range('D37').Value = "'2.654"
$ws.Range('E37').Value = '  -1.04%  '
$ws.Range('D38').Value = '1.301.35'
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('D39').Value = "'0.01791"
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('D40').Value = "'2.733"
$ws.Range('E40').Value = '  +0.89%  '
$ws.Range('D41').Value = "'0.9279"
$ws.Range('E41').Value = '  -0.86%  '
$ws.Range('D42').Value = "'5.928"
$ws.Range('E42').Value = '  -2.72%  '
$ws.Range('D43').Value = "'0.08418"
$ws.Range('E43').Value = '  +2.28%  '
$ws.Range('D44').Value = "'0.9997"
$ws.Range('D45').Value = "'102.36"
$ws.Range('E45').Value = '  -2.08%  '
$ws.Range('D46').Value = '1.963.76'
$ws.Range('E46').Value = '  -0.64%  '
$ws.Range('E47').Value = '  -0.71%  '
$ws.Range('E48').Value = '  -1.39%  '
$ws.Range('D49').Value = "'1.749"
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('D50').Value = "'62.88"
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('D51').Value = "'0.05860"
$ws.Range('E51').Value = '  -1.16%  '

Write-Host "Updated 93 cells"
